# Updated figure for prisma flowchart
# - Bumps the auto date placeholder text (master + all layouts) from 13/10/20 to 8/7/21
# - Updates the sample-size ("n = ...") figures throughout the PRISMA flow chart on slide 1

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder ("datetimeFigureOut" field) on the slide master and on
#    every slide layout. We locate the placeholder by its PlaceholderFormat
#    Type (16 = ppPlaceholderDate) instead of assuming shape index/name.
# ---------------------------------------------------------------------------
$newDate = "8/7/21"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------------
# 2. Sample-size figures in the PRISMA flowchart on slide 1.
#    All the flowchart boxes live inside a single top-level group ("Group 1").
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$grp = $slide.Shapes.Item(1)

function Set-ParaText($shape, $paraIndex, $newText) {
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs($paraIndex, 1)
    $run = $para.Runs(1, 1)
    $run.Text = $newText
}

# TextBox 7: "Records from databases ..." / n = 1605 -> n = 1829
Set-ParaText ($grp.GroupItems.Item(5)) 3 "n = 1829"

# TextBox 8: "Records after duplicates removed" / n = 1371 -> n = 1591
Set-ParaText ($grp.GroupItems.Item(6)) 3 "n = 1591"

# TextBox 9: "Records screened" / n = 1371 -> n = 1591
Set-ParaText ($grp.GroupItems.Item(7)) 3 "n = 1591"

# TextBox 15: "Records excluded based on title and abstract (n = 1122)" -> (n = 1331)
Set-ParaText ($grp.GroupItems.Item(10)) 1 "Records excluded based on title and abstract (n = 1331)"

# TextBox 19: "Full text articles assessed for eligibility" / n = 229 -> n = 260
Set-ParaText ($grp.GroupItems.Item(12)) 3 "n = 260"

# TextBox 25: exclusion reasons list
$tb25 = $grp.GroupItems.Item(15)
Set-ParaText $tb25 1 "Full text articles excluded for following reasons (n = 160):"
Set-ParaText $tb25 2 "Not an implemented AR application (n = 65)"
Set-ParaText $tb25 3 "Not for education (n = 44)"
Set-ParaText $tb25 4 "Not interactive, collaborative or multiuser (n = 28)"

# TextBox 27: "Studies included in the literature review" / n = 92 -> n = 100
Set-ParaText ($grp.GroupItems.Item(16)) 3 "n = 100"
